# Refresh the cryptocurrency price/volume snapshot (row-by-row ranking update),
# matching the "Updated cryptos list ... with GitHub Actions" scrape commit.
#
# Source data in column D/E is stored as TEXT (inline strings) even when a
# value looks like a plain decimal (e.g. "0.663") because Excel renders
# thousand-grouped prices like "44.009.43" right next to it in the same
# column. Plain `Range.Value = "0.663"` would silently get auto-coerced to
# a real Number by Excel, flipping the cell's type — so any update whose
# new text parses as a plain number is written through Set-CellText,
# which forces the Text number format just long enough to type the value
# in as a string, then restores the default "Normal" style so no stray
# per-cell number format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $cellRef, $value) {
    $range = $ws.Range($cellRef)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Looks like a bare number (e.g. "0.663", "22.54") -- force text so
        # Excel doesn't reinterpret it as a Number cell.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        # Already unambiguous as text (URLs, names, "44.009.43", percents
        # with a trailing "%" and padding spaces, etc.)
        $range.Value = $value
    }
}

$updates = @(
    @{ Cell = 'D2'; Value = '44.009.43' },
    @{ Cell = 'E2'; Value = '  +0.23%  ' },
    @{ Cell = 'D3'; Value = '2.360.82' },
    @{ Cell = 'E3'; Value = '  +1.58%  ' },
    @{ Cell = 'E4'; Value = '  +0.09%  ' },
    @{ Cell = 'D5'; Value = '0.663' },
    @{ Cell = 'E5'; Value = '  +2.71%  ' },
    @{ Cell = 'D6'; Value = '235.26' },
    @{ Cell = 'E6'; Value = '  +0.93%  ' },
    @{ Cell = 'D7'; Value = '72.63' },
    @{ Cell = 'E7'; Value = '  +9.83%  ' },
    @{ Cell = 'E8'; Value = '  -0.01%  ' },
    @{ Cell = 'D9'; Value = '0.543' },
    @{ Cell = 'E9'; Value = '  +21.52%  ' },
    @{ Cell = 'E10'; Value = '  +1.94%  ' },
    @{ Cell = 'D11'; Value = '28.17' },
    @{ Cell = 'E11'; Value = '  +4.79%  ' },
    @{ Cell = 'D12'; Value = '2.716.75' },
    @{ Cell = 'E12'; Value = '  +2.14%  ' },
    @{ Cell = 'D13'; Value = '0.107' },
    @{ Cell = 'E13'; Value = '  +1.99%  ' },
    @{ Cell = 'D14'; Value = '16.84' },
    @{ Cell = 'E14'; Value = '  +9.71%  ' },
    @{ Cell = 'D15'; Value = '6.66' },
    @{ Cell = 'E15'; Value = '  +8.67%  ' },
    @{ Cell = 'D16'; Value = '0.882' },
    @{ Cell = 'E16'; Value = '  +5.14%  ' },
    @{ Cell = 'D17'; Value = '2.367.07' },
    @{ Cell = 'E17'; Value = '  +1.98%  ' },
    @{ Cell = 'D18'; Value = '43.907.81' },
    @{ Cell = 'E18'; Value = '  +0.27%  ' },
    @{ Cell = 'E19'; Value = '  +2.76%  ' },
    @{ Cell = 'D20'; Value = '76.36' },
    @{ Cell = 'E20'; Value = '  +3.72%  ' },
    @{ Cell = 'E21'; Value = '  +2.10%  ' },
    @{ Cell = 'D22'; Value = '251.79' },
    @{ Cell = 'E22'; Value = '  +0.93%  ' },
    @{ Cell = 'B23'; Value = 'WEMIXToken' },
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'D23'; Value = '3.80' },
    @{ Cell = 'E23'; Value = '  +0.57%  ' },
    @{ Cell = 'B24'; Value = 'Dai' },
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Cell = 'D24'; Value = '0.999' },
    @{ Cell = 'E24'; Value = '  -0.11%  ' },
    @{ Cell = 'E25'; Value = '  +1.42%  ' },
    @{ Cell = 'D26'; Value = '10.45' },
    @{ Cell = 'E26'; Value = '  +5.93%  ' },
    @{ Cell = 'D27'; Value = '2.26' },
    @{ Cell = 'E27'; Value = '  -0.60%  ' },
    @{ Cell = 'D28'; Value = '22.54' },
    @{ Cell = 'E28'; Value = '  +1.61%  ' },
    @{ Cell = 'D29'; Value = '173.22' },
    @{ Cell = 'E29'; Value = '  -0.84%  ' },
    @{ Cell = 'E30'; Value = '  +8.43%  ' },
    @{ Cell = 'E31'; Value = '  +0.26%  ' },
    @{ Cell = 'D32'; Value = '0.132' },
    @{ Cell = 'D33'; Value = '5.19' },
    @{ Cell = 'E33'; Value = '  +4.06%  ' },
    @{ Cell = 'D34'; Value = '0.0712' },
    @{ Cell = 'E34'; Value = '  +3.87%  ' },
    @{ Cell = 'D35'; Value = '5.18' },
    @{ Cell = 'E35'; Value = '  +3.89%  ' },
    @{ Cell = 'D36'; Value = '3.75' },
    @{ Cell = 'E36'; Value = '  +2.42%  ' },
    @{ Cell = 'E37'; Value = '  +2.02%  ' },
    @{ Cell = 'D38'; Value = '6.41' },
    @{ Cell = 'E38'; Value = '  -1.94%  ' },
    @{ Cell = 'D39'; Value = '0.0273' },
    @{ Cell = 'E39'; Value = '  +7.95%  ' },
    @{ Cell = 'D40'; Value = '19.28' },
    @{ Cell = 'E40'; Value = '  +11.08%  ' },
    @{ Cell = 'E41'; Value = '  -0.12%  ' },
    @{ Cell = 'E42'; Value = '  -2.29%  ' },
    @{ Cell = 'E43'; Value = '  +3.97%  ' },
    @{ Cell = 'D44'; Value = '0.0974' },
    @{ Cell = 'E44'; Value = '  +1.97%  ' },
    @{ Cell = 'E45'; Value = '  +1.87%  ' },
    @{ Cell = 'B46'; Value = 'Algorand' },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Cell = 'D46'; Value = '0.183' },
    @{ Cell = 'E46'; Value = '  +13.38%  ' },
    @{ Cell = 'B47'; Value = 'FTXToken' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' },
    @{ Cell = 'D47'; Value = '4.43' },
    @{ Cell = 'E47'; Value = '  -0.75%  ' },
    @{ Cell = 'B48'; Value = 'Aave' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Cell = 'D48'; Value = '97.83' },
    @{ Cell = 'E48'; Value = '  -0.86%  ' },
    @{ Cell = 'B49'; Value = 'Maker' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Cell = 'D49'; Value = '1.438.06' },
    @{ Cell = 'E49'; Value = '  -0.38%  ' },
    @{ Cell = 'B50'; Value = 'NEARProtocol' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Cell = 'D50'; Value = '2.31' },
    @{ Cell = 'E50'; Value = '  +0.70%  ' },
    @{ Cell = 'D51'; Value = '2.586.79' },
    @{ Cell = 'E51'; Value = '  +1.90%  ' }
)

foreach ($u in $updates) {
    Set-CellText $ws $u.Cell $u.Value
}
